$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear the "Reviewer (Staff ID)" header and its sample value, keeping cell styling.
$ws.Range("J1").ClearContents()
$ws.Range("J2").ClearContents()

# Reset the view so it is scrolled back to the top-left and the selection
# covers J1:J2 with J2 as the active cell.
$ws.Range("A1").Select()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("J1:J2").Select()
$ws.Range("J2").Activate()
